$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update selection on the sheet view
$ws.Range("H13").Select()

# New block of rows (2-6) mirroring the existing SpeedTest block, for a new
# TortureTest result set referenced by the header already in A1.
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 5

$ws.Range("A3").Value = "Win64 (low frag)"
$ws.Range("F3").Value = 515958

$ws.Range("A4").Value = "nedmalloc v1.06"
$ws.Range("F4").Value = 31089140
$ws.Range("G4").Formula = "=F4/F`$3"

$ws.Range("A6").Value = "nedmalloc v1.06 (threadcached sysalloc)"

$ws.Range("A5").Value = "nedmalloc v1.06 (patcher)"
$ws.Range("F5").Value = 30994083
$ws.Range("G5:G6").Formula = "=F5/F`$3"
$ws.Range("G6").ClearContents()

# Update the existing SpeedTest (x64) block with new benchmark numbers
$ws.Range("F12").Value = 14491780

$ws.Range("F13").Value = 37044111
$ws.Range("G13").Formula = "=F13/F`$12"

$ws.Range("A14").Value = "nedmalloc v1.06 (patcher)"
$ws.Range("F14").Value = 36643063
$ws.Range("G14").Formula = "=F14/F`$12"

$ws.Range("A15").Value = "nedmalloc v1.06 (sysalloc)"
